# Connect clients to advisors and add a few more attributes
# Inserts four new attribute columns (religion, education_school,
# education_uni, profession) between "nationality" and "age", shifting
# the existing riskgroup/question columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns at E:H (pushes old E:M -> I:Q), carrying the
# existing formulas/formatting along for the ride and auto-adjusting
# the shared-formula references (SUM(I2:M2) -> SUM(M2:Q2), etc.).
$ws.Range("E1:H1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("E1").Value = "religion"
$ws.Range("F1").Value = "education_school"
$ws.Range("G1").Value = "education_uni"
$ws.Range("H1").Value = "profession"

# Match the saved selection/active cell from the edited workbook.
$ws.Range("H2").Select() | Out-Null
